# Update 江西-漫展信息.xlsx "展览" and "全部类型" sheets:
# - Remove the canceled "江西·JMG（江西广电）第二届UP动漫游戏博览会-钱琛专场见面会"
#   event (row with A-index 4, B=2024-10-04) from both sheets; every
#   following row shifts up by one row.
# - Refresh the "想去人数" (column F) counters with the newer scrape values
#   for the rows that shifted up.
# - Renumber the leading index column (A) so it stays 0,1,2,3,... after the
#   deletion.

$wb = $excel.ActiveWorkbook

# Sheets that contain the removed row + shifted rows, and the 1-based row
# (after the deletion has already happened) -> new "想去人数" value.
$sheetUpdates = @{
    "展览"   = @{5=36; 7=150; 8=629; 10=194; 11=1320; 12=24; 13=2936; 14=446; 15=645}
    "全部类型" = @{5=36; 7=150; 8=629; 11=194; 12=1320; 13=24; 14=2936; 15=446; 16=645}
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets($sheetName)

    # Rows 2 and 3 get refreshed "想去人数" values regardless of the deletion.
    $ws.Cells.Item(2, 6).Value = 4502
    $ws.Cells.Item(3, 6).Value = 855

    # Delete the canceled event's row (row 5: 2024-10-04 钱琛专场见面会).
    # This shifts every following row up by one and shrinks the sheet's
    # used range automatically.
    $ws.Rows(5).Delete()

    # Fix up the index column (A) so it keeps counting 0,1,2,3,... after
    # the shift (row deletion does not renumber plain literal values).
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 5; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply the refreshed "想去人数" (column F) values on the rows that
    # shifted up, per the newer scrape.
    $updates = $sheetUpdates[$sheetName]
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
